# "added pressures and temps"
#
# The Translation sheet (table "Table8"/table3.xml) lists gauge text entries
# as (Text ID, Typography Name, Alignment, GB, Direction) rows. This edit
# renumbers a block of auto-generated "SingleUseIdNN" Text IDs, re-pairs a
# couple of typography rows, and appends four new descriptor/value rows for
# FUEL T / BAT / MIL / OIL pressures & temperatures.
#
# NOTE: several of the "GB" (column E) values are the literal text "0" or
# "0.0" (placeholders for numeric gauge values). Assigning those through
# Range.Value would be auto-parsed as numbers (losing the "0.0" text and
# picking up a new NumberFormat style), so those specific cells are filled
# by copy/pasting from an existing text cell that already holds the exact
# same string - this carries the string over verbatim with no style churn.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 60 currently holds the only "0.0" text cell (E60); row 5 always holds a
# plain text "0" (E5, untouched by this edit). Stash "0.0" into a scratch
# cell first since E60 itself is about to be overwritten with "0".
$ws.Range("E60").Copy($ws.Range("Z1000"))

# --- Existing-row edits -----------------------------------------------------

# Rows 34 / 44: "0" value placeholder -> "0.0"
$ws.Range("Z1000").Copy($ws.Range("E34"))
$ws.Range("Z1000").Copy($ws.Range("E44"))

# Row 49: same "0" -> "0.0" tweak
$ws.Range("Z1000").Copy($ws.Range("E49"))

# Rows 50-53: Text IDs renumbered (SingleUseId78-81 -> 80,81,82,83)
$ws.Cells.Item(50, 2).Value = "SingleUseId80"
$ws.Cells.Item(51, 2).Value = "SingleUseId81"
$ws.Cells.Item(52, 2).Value = "SingleUseId82"
$ws.Cells.Item(53, 2).Value = "SingleUseId83"

# Rows 54-57: the "Default"/"small" typography pairing is swapped & renumbered
$ws.Cells.Item(54, 2).Value = "SingleUseId84"
$ws.Cells.Item(54, 3).Value = "small"
$ws.Cells.Item(54, 5).Value = "<value> "

$ws.Cells.Item(55, 2).Value = "SingleUseId85"
$ws.Cells.Item(55, 3).Value = "small"

$ws.Cells.Item(56, 2).Value = "SingleUseId88"
$ws.Cells.Item(56, 3).Value = "Default"
$ws.Cells.Item(56, 5).Value = "<value>"

$ws.Cells.Item(57, 2).Value = "SingleUseId89"
$ws.Cells.Item(57, 3).Value = "Default"

# Row 58: becomes the "Large" / Left / "0.0" value row
$ws.Cells.Item(58, 2).Value = "SingleUseId90"
$ws.Cells.Item(58, 3).Value = "Large"
$ws.Cells.Item(58, 4).Value = "Left"
$ws.Range("Z1000").Copy($ws.Range("E58"))

# Row 59: becomes a new "RpmText" / Center / "<value>" descriptor row
$ws.Cells.Item(59, 2).Value = "SingleUseId91"
$ws.Cells.Item(59, 3).Value = "RpmText"
$ws.Cells.Item(59, 4).Value = "Center"
$ws.Cells.Item(59, 5).Value = "<value>"

# Row 60: becomes the matching "RpmText" / Left / "0" value row
$ws.Cells.Item(60, 2).Value = "SingleUseId92"
$ws.Cells.Item(60, 3).Value = "RpmText"
$ws.Range("E5").Copy($ws.Range("E60"))

# Scratch cell no longer needed
$ws.Range("Z1000").ClearContents()

# --- New rows: pressures and temps ------------------------------------------

$ws.Cells.Item(61, 2).Value = "SingleUseId93"
$ws.Cells.Item(61, 3).Value = "Default"
$ws.Cells.Item(61, 4).Value = "Left"
$ws.Cells.Item(61, 5).Value = "FUEL T"
$ws.Cells.Item(61, 6).Value = "LTR"

$ws.Cells.Item(62, 2).Value = "SingleUseId94"
$ws.Cells.Item(62, 3).Value = "Values"
$ws.Cells.Item(62, 4).Value = "Left"
$ws.Cells.Item(62, 5).Value = "BAT"
$ws.Cells.Item(62, 6).Value = "LTR"

$ws.Cells.Item(63, 2).Value = "SingleUseId95"
$ws.Cells.Item(63, 3).Value = "Values"
$ws.Cells.Item(63, 4).Value = "Left"
$ws.Cells.Item(63, 5).Value = "MIL"
$ws.Cells.Item(63, 6).Value = "LTR"

$ws.Cells.Item(64, 2).Value = "SingleUseId96"
$ws.Cells.Item(64, 3).Value = "Values"
$ws.Cells.Item(64, 4).Value = "Left"
$ws.Cells.Item(64, 5).Value = "OIL"
$ws.Cells.Item(64, 6).Value = "LTR"
